$wb = $excel.ActiveWorkbook

# --- "About" sheet: update/insert explanatory notes about the India
# pumped-hydro (repurposed geothermal) flexibility points correction ---
$ws1 = $wb.Worksheets.Item("About")

# Row 33 text tweaked (trailing period removed) and now continues onto two
# new lines (34 and 35) explaining the 0.5 flexibility value for geothermal.
$ws1.Range("A33").Value = "For India, Flag zero for hydro because we track pumped hydro separately"
$ws1.Range("A35").Value = "the flexibility points as peaker plants, we use a value of 0.5 for geothermal here."
$ws1.Range("A34").Value = "as the geothermal plant type. Because pumped hydro plants provide half (see elec/FPC)"

# Row 36 becomes the blank spacer row (content that used to live there moves
# down to rows 37-40).
$ws1.Range("A36").ClearContents()
$ws1.Range("A37").Value = "Natural gas nonpeakers (CCGTs) are flagged for flexibility because"
$ws1.Range("A38").Value = "in India they are are used for peaking as well as a low level of baseload."
$ws1.Range("A39").Value = "Petroleum/diesel plants are flagged as 1 as they are used for balancing "
$ws1.Range("A40").Value = "due to fast ramping nature."

# --- "BPaFF-BDTPTPF" sheet: geothermal (repurposed pumped hydro) now
# provides only half (0.5) flexibility points instead of a full point ---
$ws3 = $wb.Worksheets.Item("BPaFF-BDTPTPF")
$ws3.Range("B10").Value = 0.5
